# Apply the updated crypto symbol list values (Price + Volume(1h) columns).
# Target cells are re-written as literal text, matching the source feed's
# string formatting (e.g. trailing zeros in prices, "%" suffixed deltas),
# while keeping each cell's original (default/unstyled) formatting intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $CellRef, $Text) {
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

Set-TextValue $ws "D2" "307.60"
Set-TextValue $ws "E2" "-4.17%"
Set-TextValue $ws "D3" "39.87"
Set-TextValue $ws "E3" "-6.61%"
Set-TextValue $ws "D4" "5.113"
Set-TextValue $ws "E4" "-1.72%"
Set-TextValue $ws "D5" "0.07717"
Set-TextValue $ws "E5" "-5.60%"
Set-TextValue $ws "D6" "4.231"
Set-TextValue $ws "E6" "-1.84%"
Set-TextValue $ws "D7" "1.606"
Set-TextValue $ws "E7" "-11.43%"
Set-TextValue $ws "D8" "0.8938"
Set-TextValue $ws "E8" "-4.30%"
Set-TextValue $ws "D9" "0.1005"
Set-TextValue $ws "E9" "-9.35%"
Set-TextValue $ws "D10" "0.1734"
Set-TextValue $ws "E10" "-6.77%"
Set-TextValue $ws "D11" "0.09024"
Set-TextValue $ws "E11" "-3.76%"
Set-TextValue $ws "D12" "0.04451"
Set-TextValue $ws "E12" "-5.58%"
Set-TextValue $ws "D14" "0.001269"
Set-TextValue $ws "E14" "-1.59%"
Set-TextValue $ws "D15" "0.005806"
Set-TextValue $ws "E15" "0.04%"
Set-TextValue $ws "D19" "0.3317"
Set-TextValue $ws "E19" "-0.93%"
Set-TextValue $ws "D20" "7.050"
Set-TextValue $ws "E20" "-5.11%"
Set-TextValue $ws "D23" "0.04142"
Set-TextValue $ws "E23" "-0.06%"
Set-TextValue $ws "D24" "0.001208"
Set-TextValue $ws "E24" "-2.92%"
Set-TextValue $ws "D25" "0.004060"
Set-TextValue $ws "E25" "-5.55%"
Set-TextValue $ws "D26" "0.0001302"
Set-TextValue $ws "E26" "8.43%"
Set-TextValue $ws "D38" "0.02341"
Set-TextValue $ws "E38" "-13.10%"
Set-TextValue $ws "D39" "0.05197"
Set-TextValue $ws "E39" "-6.26%"
Set-TextValue $ws "D40" "0.007917"
Set-TextValue $ws "E40" "-2.61%"
Set-TextValue $ws "D42" "0.006249"
Set-TextValue $ws "E42" "-4.50%"
Set-TextValue $ws "D43" "0.001953"
Set-TextValue $ws "E43" "-6.43%"
Set-TextValue $ws "D44" "0.008226"
Set-TextValue $ws "E44" "-0.34%"
Set-TextValue $ws "D45" "0.3329"
Set-TextValue $ws "E45" "-4.67%"
Set-TextValue $ws "D46" "0.00006513"
Set-TextValue $ws "E46" "-5.99%"
Set-TextValue $ws "D49" "0.003445"
Set-TextValue $ws "E49" "2.91%"
Set-TextValue $ws "D50" "0.00002104"
Set-TextValue $ws "E50" "0.09%"
Set-TextValue $ws "D51" "0.0002004"
Set-TextValue $ws "E51" "0.09%"
Set-TextValue $ws "E16" "2,411.83%"
Set-TextValue $ws "E17" "-0.08%"
Set-TextValue $ws "E18" "-3.53%"
Set-TextValue $ws "E21" "-2.40%"
Set-TextValue $ws "E22" "8.30%"
Set-TextValue $ws "E41" "-5.58%"
Set-TextValue $ws "E47" "0.09%"
Set-TextValue $ws "E48" "98.24%"
